$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (66, 67, 68) to the list
$ws.Range("A66").Value = "Dorsaf"
$ws.Range("B66").Value = "Sallami"
$ws.Range("C66").Value = "Université de Montréal"
$ws.Range("D66").Value = "Canada"
$ws.Range("E66").Value = "wUa3IWgAAAAJ"
$ws.Range("F66").Value = "F"
$ws.Range("G66").Value = 1996
$ws.Range("H66").Value = "Informatique, Mathématiques et Ingénierie"

$ws.Range("A67").Value = "Ilef"
$ws.Range("B67").Value = "Romdhani"
$ws.Range("C67").Value = "Université de Monastir"
$ws.Range("D67").Value = "Tunisie"
$ws.Range("E67").Value = "_ynJvKAAAAAJ"
$ws.Range("F67").Value = "F"
$ws.Range("G67").Value = 1998
$ws.Range("H67").Value = "Médecine, Biologie et Sciences de la Santé"

$ws.Range("A68").Value = "Ilyes"
$ws.Range("B68").Value = "Rezgui"
$ws.Range("C68").Value = "Université de Tunis El Manar"
$ws.Range("D68").Value = "Tunisie"
$ws.Range("E68").Value = "2RI_ZRkAAAAJ"
$ws.Range("F68").Value = "M"
$ws.Range("G68").Value = 2001
$ws.Range("H68").Value = "Informatique, Mathématiques et Ingénierie"

# Copy the font style (Arial 8pt) from existing Genre column cells to the new ones
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F66:F68").PasteSpecial(-4122) | Out-Null

# Restore the view settings to match the new scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 51
$ws.Range("D66").Select() | Out-Null
